$wb = $excel.ActiveWorkbook

$sheetPedInfluenza = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$sheetAdultInfluenza = $wb.Worksheets.Item("Adult Influenza Vaccine ")

# Pediatric Influenza Vaccine sheet: collapse multi-line strings to single-line (newline -> space)
$sheetPedInfluenza.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$sheetPedInfluenza.Range("B6").Value = "Fluarix Preservative-Free"
$sheetPedInfluenza.Range("B9").Value = "FluMist No Preservative"
$sheetPedInfluenza.Range("B10").Value = "Afluria No Preservative"
$sheetPedInfluenza.Range("H10").Value = "Merck (CSL product)"

# Adult Influenza Vaccine sheet: collapse multi-line strings to single-line (newline -> space)
$sheetAdultInfluenza.Range("B5").Value = "Agriflu No Preservative"
$sheetAdultInfluenza.Range("B7").Value = "Fluvirin Preservative-free"
$sheetAdultInfluenza.Range("B8").Value = "Fluarix Preservative-free"
$sheetAdultInfluenza.Range("B10").Value = "Flumist No Preservative"
